# Update cryptocurrency price/volume data per upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text numeric-looking strings (e.g. '42.612.77',
# '1.00', '0.0968') that must stay text, not be coerced to numbers, so force the
# Text number format on those cells before writing the new values.
$priceCells = @(
    "D2", "D3", "D5", "D7", "D9", "D10", "D11", "D12",
    "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21",
    "D22", "D23", "D24", "D26", "D27", "D31", "D32", "D33",
    "D35", "D36", "D38", "D40", "D41", "D42", "D43", "D45",
    "D46", "D49", "D50", "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values (Coin name / Link / Price / Volume(1h)).
$ws.Range("D2").Value = '42.612.77'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '2.286.24'
$ws.Range("E3").Value = '  +4.44%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '250.54'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("E6").Value = '  +2.48%  '
$ws.Range("D7").Value = '72.15'
$ws.Range("E7").Value = '  +7.43%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.645'
$ws.Range("E9").Value = '  +3.68%  '
$ws.Range("D10").Value = '38.82'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").Value = '0.0968'
$ws.Range("E11").Value = '  +4.01%  '
$ws.Range("D12").Value = '59.07'
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Value = '7.33'
$ws.Range("E13").Value = '  +4.99%  '
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").Value = '2.627.99'
$ws.Range("E15").Value = '  +4.35%  '
$ws.Range("D16").Value = '15.09'
$ws.Range("E16").Value = '  +4.32%  '
$ws.Range("D17").Value = '0.881'
$ws.Range("E17").Value = '  +2.81%  '
$ws.Range("D18").Value = '2.285.20'
$ws.Range("E18").Value = '  +3.86%  '
$ws.Range("D19").Value = '42.556.77'
$ws.Range("E19").Value = '  +2.84%  '
$ws.Range("E20").Value = '  +4.72%  '
$ws.Range("D21").Value = '6.30'
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("D22").Value = '72.43'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").Value = '2.28'
$ws.Range("E23").Value = '  +12.69%  '
$ws.Range("D24").Value = '235.18'
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").Value = '11.63'
$ws.Range("E26").Value = '  +3.29%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("E30").Value = '  +6.69%  '
$ws.Range("D31").Value = '167.61'
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("D32").Value = '21.10'
$ws.Range("E32").Value = '  +4.33%  '
$ws.Range("D33").Value = '6.44'
$ws.Range("E33").Value = '  +10.54%  '
$ws.Range("E34").Value = '  +6.17%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.0804'
$ws.Range("E35").Value = '  +2.86%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '31.56'
$ws.Range("E36").Value = '  +22.86%  '
$ws.Range("E37").Value = '  +3.29%  '
$ws.Range("D38").Value = '4.74'
$ws.Range("E38").Value = '  +14.53%  '
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("D40").Value = '0.0307'
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").Value = '14.03'
$ws.Range("E41").Value = '  +15.28%  '
$ws.Range("D42").Value = '2.32'
$ws.Range("E42").Value = '  +5.44%  '
$ws.Range("D43").Value = '5.96'
$ws.Range("E43").Value = '  +6.86%  '
$ws.Range("E44").Value = '  +9.48%  '
$ws.Range("D45").Value = '9.22'
$ws.Range("E45").Value = '  +8.34%  '
$ws.Range("D46").Value = '61.97'
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("E47").Value = '  -3.80%  '
$ws.Range("E48").Value = '  +3.64%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '1.18'
$ws.Range("E49").Value = '  +3.09%  '
$ws.Range("B50").Value = 'BinanceUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").Value = '96.84'
$ws.Range("E51").Value = '  +6.33%  '
